$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$reply1 = "Beste afzender,`nDank voor je bericht. Onze kantooruren zijn van maandag t/m vrijdag van 9:00-17:00 uur. Buiten deze tijden streven we naar spoedige reactie. Voor verdere vragen, neem gerust contact met ons op.`nMet vriendelijke groet,`n[E-mailassistent]"
$reply2 = "Beste [Naam],`nDank voor je bericht. Onze kantooruren zijn van maandag t/m vrijdag van 9:00-17:00 uur. Buiten deze tijden streven we naar een spoedige reactie. Voor verdere vragen, neem gerust contact met ons op.`nMet vriendelijke groet,`n[E-mailassistent]"
$reply3 = "Beste [Naam],`nDank voor je bericht. Onze kantooruren zijn van maandag t/m vrijdag van 9:00-17:00 uur. Buiten deze tijden streven we naar een spoedige reactie. Voor verdere vragen, neem gerust contact met ons op.`nMet vriendelijke groet,`nE-mailassistent"

$rows = @(
    @{ r=23; A="Re: Re: Wat zijn jullie openingstijden?"; C=$reply1; D="Informatieaanvraag"; E=$reply1; F="2025-06-17 20:44:19"; G="Ja" },
    @{ r=24; A="Re: Re: Re: Wat zijn jullie openingstijden?"; C=$reply1; D="Informatieaanvraag"; E=$reply2; F="2025-06-17 20:44:34"; G="Ja" },
    @{ r=25; A="Re: Re: Re: Re: Wat zijn jullie openingstijden?"; C=$reply2; D="Informatieaanvraag"; E=$reply3; F="2025-06-17 20:46:24"; G="Ja" },
    @{ r=26; A="Re: Re: Re: Re: Re: Wat zijn jullie openingstijden?"; C=$reply3; D="Informatieaanvraag"; E=$reply3; F="2025-06-17 20:46:41"; G="Ja" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = "mailmind.test@zohomail.eu"
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}

$dash.Range("B2").Value = 13

foreach ($fc in $ws.Range("D2:D22").FormatConditions) {
    $fc.ModifyAppliesToRange($ws.Range("D2:D26"))
}
foreach ($fc in $ws.Range("G2:G22").FormatConditions) {
    $fc.ModifyAppliesToRange($ws.Range("G2:G26"))
}
